$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 116
$ws.Range("H116").Value = 4230
$ws.Range("I116").Value = 2528.5715
$ws.Range("J116").Value = 5146.154
$ws.Range("K116").Value = 2528.5715
$ws.Range("L116").Value = 5146.154
$ws.Range("M116").Value = 913.4285
$ws.Range("N116").Value = -12030.154

# Row 133
$ws.Range("H133").Value = 32977.5
$ws.Range("J133").Value = 32977.5
$ws.Range("L133").Value = 32977.5
$ws.Range("N133").Value = -43097.5

# Row 138
$ws.Range("H138").Value = 2787.7236
$ws.Range("I138").Value = 2460.68
$ws.Range("J138").Value = 2948.0393
$ws.Range("K138").Value = 7382.039999999999
$ws.Range("L138").Value = 8844.117899999999
$ws.Range("M138").Value = -2242.039999999999
$ws.Range("N138").Value = -19124.1179

# Row 141
$ws.Range("H141").Value = 3177.8262
$ws.Range("I141").Value = 1274.7059
$ws.Range("J141").Value = 8570
$ws.Range("K141").Value = 3824.1177
$ws.Range("L141").Value = 25710
$ws.Range("M141").Value = 1355.8823
$ws.Range("N141").Value = -36070

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 9097.5
$ws.Range("I32").Value = 8339.954
$ws.Range("J32").Value = 14097.3
$ws.Range("K32").Value = 8339.954
$ws.Range("L32").Value = 14097.3
$ws.Range("M32").Value = -8052.954
$ws.Range("N32").Value = -14671.3

# Row 63
$ws.Range("H63").Value = 2870
$ws.Range("I63").Value = 2360
$ws.Range("J63").Value = 4400
$ws.Range("K63").Value = 2360
$ws.Range("L63").Value = 4400
$ws.Range("M63").Value = -1674
$ws.Range("N63").Value = -5772

# Row 66
$ws.Range("H66").Value = 2870
$ws.Range("I66").Value = 2360
$ws.Range("J66").Value = 4400
$ws.Range("K66").Value = 11800
$ws.Range("L66").Value = 22000
$ws.Range("M66").Value = -8368
$ws.Range("N66").Value = -28864

# Row 92
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("N92").Value = 0
$ws.Range("L92").ClearContents()

# Row 132
$ws.Range("H132").Value = 7354264
$ws.Range("I132").Value = 10204987
$ws.Range("J132").Value = 2400.842
$ws.Range("K132").Value = 30614961
$ws.Range("L132").Value = 7202.526
$ws.Range("M132").Value = -30612431
$ws.Range("N132").Value = -12262.526

# Row 138
$ws.Range("H138").Value = 43000
$ws.Range("J138").Value = 43000
$ws.Range("L138").Value = 43000
$ws.Range("N138").Value = -53280

# Row 139
$ws.Range("H139").Value = 47389
$ws.Range("J139").Value = 47389
$ws.Range("L139").Value = 47389
$ws.Range("N139").Value = -57669

$ws = $wb.Worksheets.Item("BSM")
# Row 35
$ws.Range("H35").Value = 34990
$ws.Range("J35").Value = 34990
$ws.Range("L35").Value = 34990
$ws.Range("N35").Value = -35610

# Row 70
$ws.Range("H70").Value = 136666.67
$ws.Range("I70").Value = 130000
$ws.Range("J70").Value = 150000
$ws.Range("K70").Value = 130000
$ws.Range("L70").Value = 150000
$ws.Range("M70").Value = -129707
$ws.Range("N70").Value = -150586

# Row 73
$ws.Range("H73").Value = 136666.67
$ws.Range("I73").Value = 130000
$ws.Range("J73").Value = 150000
$ws.Range("K73").Value = 130000
$ws.Range("L73").Value = 150000
$ws.Range("M73").Value = -128986
$ws.Range("N73").Value = -152028

# Row 82
$ws.Range("H82").Value = 7916.6665
$ws.Range("I82").Value = 1875
$ws.Range("J82").Value = 20000
$ws.Range("K82").Value = 1875
$ws.Range("L82").Value = 20000
$ws.Range("M82").Value = -1492
$ws.Range("N82").Value = -20766

# Row 85
$ws.Range("H85").Value = 7916.6665
$ws.Range("I85").Value = 1875
$ws.Range("J85").Value = 20000
$ws.Range("K85").Value = 1875
$ws.Range("L85").Value = 20000
$ws.Range("M85").Value = -549
$ws.Range("N85").Value = -22652

# Row 92
$ws.Range("H92").Value = 42569.5
$ws.Range("J92").Value = 42569.5
$ws.Range("L92").Value = 42569.5
$ws.Range("N92").Value = -47561.5

# Row 105
$ws.Range("H105").Value = 2273.55
$ws.Range("I105").Value = 2045.5555
$ws.Range("J105").Value = 2460.0908
$ws.Range("K105").Value = 2045.5555
$ws.Range("L105").Value = 2460.0908
$ws.Range("M105").Value = -298.5554999999999
$ws.Range("N105").Value = -5954.0908

# Row 112
$ws.Range("H112").Value = 46661
$ws.Range("J112").Value = 46661
$ws.Range("L112").Value = 46661
$ws.Range("N112").Value = -49615

# Row 132
$ws.Range("H132").Value = 42040
$ws.Range("J132").Value = 42040
$ws.Range("L132").Value = 42040
$ws.Range("N132").Value = -52160

# Row 133
$ws.Range("H133").Value = 48249.75
$ws.Range("J133").Value = 48249.75
$ws.Range("L133").Value = 48249.75
$ws.Range("N133").Value = -58369.75

$ws = $wb.Worksheets.Item("CRP")
# Row 20
$ws.Range("H20").Value = 49925
$ws.Range("J20").Value = 49925
$ws.Range("L20").Value = 49925
$ws.Range("N20").Value = -50397

# Row 30
$ws.Range("H30").Value = 49925
$ws.Range("J30").Value = 49925
$ws.Range("L30").Value = 49925
$ws.Range("N30").Value = -50107

# Row 128
$ws.Range("H128").Value = 49925
$ws.Range("J128").Value = 49925
$ws.Range("L128").Value = 49925
$ws.Range("N128").Value = -59885

# Row 137
$ws.Range("H137").Value = 45796.668
$ws.Range("J137").Value = 45796.668
$ws.Range("L137").Value = 45796.668
$ws.Range("N137").Value = -55996.668

$ws = $wb.Worksheets.Item("CUL")
# Row 82
$ws.Range("H82").Value = 3980
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 3980
$ws.Range("K82").Value = 0
$ws.Range("M82").Value = 11940
$ws.Range("N82").Value = -12752
$ws.Range("L82").ClearContents()

# Row 85
$ws.Range("H85").Value = 3980
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 3980
$ws.Range("K85").Value = 0
$ws.Range("M85").Value = 11940
$ws.Range("N85").Value = -14748
$ws.Range("L85").ClearContents()

# Row 88
$ws.Range("H88").Value = 6885.8887
$ws.Range("J88").Value = 6885.8887
$ws.Range("L88").Value = 20657.6661
$ws.Range("N88").Value = -21513.6661

# Row 91
$ws.Range("H91").Value = 6885.8887
$ws.Range("J91").Value = 6885.8887
$ws.Range("L91").Value = 20657.6661
$ws.Range("N91").Value = -23621.6661

# Row 94
$ws.Range("H94").Value = 2757.0715
$ws.Range("I94").Value = 1300
$ws.Range("J94").Value = 2999.9167
$ws.Range("K94").Value = 3900
$ws.Range("L94").Value = 8999.750100000001
$ws.Range("M94").Value = -3224
$ws.Range("N94").Value = -10351.7501

# Row 96
$ws.Range("H96").Value = 4500
$ws.Range("J96").Value = 4500
$ws.Range("L96").Value = 13500
$ws.Range("N96").Value = -17618

# Row 97
$ws.Range("H97").Value = 894.25
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 894.25
$ws.Range("K97").Value = 0
$ws.Range("M97").Value = 2682.75
$ws.Range("N97").Value = -3674.75
$ws.Range("L97").ClearContents()

# Row 103
$ws.Range("H103").Value = 1404.1666
$ws.Range("I103").Value = 806.25
$ws.Range("J103").Value = 2600
$ws.Range("K103").Value = 2418.75
$ws.Range("L103").Value = 7800
$ws.Range("M103").Value = -1539.75
$ws.Range("N103").Value = -9558

# Row 131
$ws.Range("H131").Value = 2819.6035
$ws.Range("I131").Value = 7519.2144
$ws.Range("J131").Value = 1324.2727
$ws.Range("K131").Value = 22557.6432
$ws.Range("L131").Value = 3972.8181
$ws.Range("M131").Value = -17517.6432
$ws.Range("N131").Value = -14052.8181

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 1704.4615
$ws.Range("I132").Value = 1144.2941
$ws.Range("K132").Value = 3432.8823
$ws.Range("M132").Value = -902.8823000000002

# Row 135
$ws.Range("H135").Value = 41811.43
$ws.Range("J135").Value = 41811.43
$ws.Range("L135").Value = 41811.43
$ws.Range("N135").Value = -51951.43

# Row 136
$ws.Range("H136").Value = 85163
$ws.Range("J136").Value = 85163
$ws.Range("L136").Value = 255489
$ws.Range("N136").Value = -260589

# Row 138
$ws.Range("H138").Value = 53500
$ws.Range("J138").Value = 53500
$ws.Range("L138").Value = 53500
$ws.Range("N138").Value = -63780

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 1983.8032
$ws.Range("I132").Value = 1219.0714
$ws.Range("J132").Value = 3674.2632
$ws.Range("K132").Value = 3657.2142
$ws.Range("L132").Value = 11022.7896
$ws.Range("M132").Value = -1127.2142
$ws.Range("N132").Value = -16082.7896

# Row 134
$ws.Range("H134").Value = 47747.75
$ws.Range("J134").Value = 47747.75
$ws.Range("L134").Value = 47747.75
$ws.Range("N134").Value = -57887.75

# Row 136
$ws.Range("H136").Value = 1585.2894
$ws.Range("I136").Value = 1352.8889
$ws.Range("K136").Value = 4058.6667
$ws.Range("M136").Value = -1508.6667

$ws = $wb.Worksheets.Item("WVR")
# Row 135
$ws.Range("H135").Value = 39900
$ws.Range("J135").Value = 39900
$ws.Range("L135").Value = 39900
$ws.Range("N135").Value = -50040

# Row 137
$ws.Range("H137").Value = 52999
$ws.Range("J137").Value = 52999
$ws.Range("L137").Value = 52999
$ws.Range("N137").Value = -63199

# Row 138
$ws.Range("H138").Value = 44443.332
$ws.Range("J138").Value = 44443.332
$ws.Range("L138").Value = 44443.332
$ws.Range("N138").Value = -54723.332
